$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Delete the "20200916_Morning_Noise" sheet
$wb.Worksheets.Item("20200916_Morning_Noise").Delete()

# Rename "20201207_alpha_sim" to "simulations"
$wb.Worksheets.Item("20201207_alpha_sim").Name = "simulations"

# Make "simulations" the active sheet/tab
$wb.Worksheets.Item("simulations").Activate()
